# Update annotations for Ruilin
# 1) Row 96's politeness_score (B96) was stored as text "3"; correct it to a
#    genuine numeric value 3 (matches the target OOXML: t="n" / <v>3</v>).
# 2) Append a new annotation row (97) with the new review comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B96: make it a real number instead of a text "3" ---
$ws.Cells.Item(96, 2).Value = 3

# --- Append new row 97 ---
$ws.Cells.Item(97, 1).Value = "Ruilin"

# politeness_score (B97) must stay a text "3" (like the original B96 was),
# so force text formatting before assigning, then restore the default style
# so no extra formatting is left on the cell.
$ws.Cells.Item(97, 2).NumberFormat = "@"
$ws.Cells.Item(97, 2).Value = "3"
$ws.Cells.Item(97, 2).Style = "Normal"

$ws.Cells.Item(97, 3).Value = "无"
$ws.Cells.Item(97, 4).Value = "FBK"
$ws.Cells.Item(97, 5).Value = "WRI"
$ws.Cells.Item(97, 6).Value = "3419a239-823d-4d38-8055-389a9317394a"
$ws.Cells.Item(97, 7).Value = "SJa9iHgAZ_annotated.xlsx"
$ws.Cells.Item(97, 8).Value = "To address Reviewer 2 comment on iterative inference in shared Resnet, we added two sections in Appendix reporting metrics (cosine loss, accuracy, l1 ratio) on shared Resnet, and on the unrolled to more steps Resnet."
